$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.046.23"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +3.65%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.804.57"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +4.27%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9984"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.85"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9986"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5497"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +12.43%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3799"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +8.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.10"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07593"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +4.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.134"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +7.98%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9987"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.17"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +5.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.230"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +5.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.800.14"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +4.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.153"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.02"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +5.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001081"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +3.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06501"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9984"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.17"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.993"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +5.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.049.88"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.25"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.86%  "
$ws.Range("E25").Value = "  +0.79%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.61"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.85%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.04"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.387"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +14.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.009.09"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +4.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.06"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.150"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +8.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1038"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +11.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.759"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +6.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.599"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02301"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +5.00%  "
$ws.Range("B36").Value = "Algorand"
$ws.Range("C36").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2124"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +6.12%  "
$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.656"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +15.22%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "11.52"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +4.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.022"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +4.83%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06041"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6304"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +4.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9980"
$ws.Range("D42").ClearFormats()
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.399"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.152"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +4.69%  "
$ws.Range("E45").Value = "  +4.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5925"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +4.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.670"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "121.89"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.927"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.64%  "
$ws.Range("E50").Value = "  +2.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06789"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.19%  "
